# Data-wrangling transform: pivot the old 4-column x 2-data-row table
# (Year / 385,193.. / 235,319..) into a 2-column x 4-data-row table with
# an "Unnamed" header row (pandas-style export), values unchanged.
#
# Old layout:
#   A1=Year
#   A2=385,193  B2=83   C2=2,727,635  D2=586
#   A3=(blank)
#   A4=235,319  B4=50   C4=2,631,789  D4=558
#   A5=(blank)
#
# New layout:
#   A1=Unnamed: 0  B1=2019  C1=Unnamed: 1  D1=2018  E1=Unnamed: 2
#   B2=385,193     D2=235,319
#   B3=83          D3=50
#   B4=2,727,635   D4=2,631,789
#   B5=586         D5=558

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163
$xlPasteFormats = -4122

# --- Stage the two existing data rows off to the side (row-shaped ranges,
# matching the source shape) so the in-place column/row permutation below
# can't clobber a source cell before it has been read. ---
$ws.Range("A2:D2").Copy($ws.Range("G1:J1"))
$ws.Range("A4:D4").Copy($ws.Range("G2:J2"))

# --- Header row (row 1). A1 already carries the bold/boxed style, so only
# its text changes; B1/C1/D1/E1 pick the same style up via a format-only
# paste (keeps styles.xml free of ad-hoc new styles). ---
$ws.Range("A1").Copy()
$ws.Range("B1:E1").PasteSpecial($xlPasteFormats)

$ws.Range("A1").Value = "Unnamed: 0"
$ws.Range("C1").Value = "Unnamed: 1"
$ws.Range("E1").Value = "Unnamed: 2"
# "2019"/"2018" look numeric, so force text (as in the source data) with a
# leading apostrophe - matches the shared-string <t> entries in the diff.
$ws.Range("B1").Value = "'2019"
$ws.Range("D1").Value = "'2018"

# --- Column B (rows 2-5) <- old row 2 (A2,B2,C2,D2), values only (no
# style) so the cells pick up the workbook's default formatting. ---
$ws.Range("G1").Copy()
$ws.Range("B2").PasteSpecial($xlPasteValues)
$ws.Range("H1").Copy()
$ws.Range("B3").PasteSpecial($xlPasteValues)
$ws.Range("I1").Copy()
$ws.Range("B4").PasteSpecial($xlPasteValues)
$ws.Range("J1").Copy()
$ws.Range("B5").PasteSpecial($xlPasteValues)

# --- Column D (rows 2-5) <- old row 4 (A4,B4,C4,D4), values only. ---
$ws.Range("G2").Copy()
$ws.Range("D2").PasteSpecial($xlPasteValues)
$ws.Range("H2").Copy()
$ws.Range("D3").PasteSpecial($xlPasteValues)
$ws.Range("I2").Copy()
$ws.Range("D4").PasteSpecial($xlPasteValues)
$ws.Range("J2").Copy()
$ws.Range("D5").PasteSpecial($xlPasteValues)

# --- Drop cells that are no longer part of the table: old A2:A5 leftover
# styling, old C2/C4 values, and the staging area. ---
$ws.Range("A2:A5").Clear()
$ws.Range("C2").Clear()
$ws.Range("C4").Clear()
$ws.Range("G1:J2").Clear()
